$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1788
$ws.Range("F3").Value = 10472
$ws.Range("F6").Value = 624
$ws.Range("F8").Value = 1761
$ws.Range("F12").Value = 568
$ws.Range("F15").Value = 39
$ws.Range("F16").Value = 1042
$ws.Range("F18").Value = 133
$ws.Range("F19").Value = 465
$ws.Range("F20").Value = 465
$ws.Range("F22").Value = 379
$ws.Range("F23").Value = 93
$ws.Range("F24").Value = 1160
$ws.Range("F25").Value = 1138
$ws.Range("F26").Value = 1270
$ws.Range("F27").Value = 230
$ws.Range("F28").Value = 1454
$ws.Range("F29").Value = 737
$ws.Range("F30").Value = 272
$ws.Range("F31").Value = 37
$ws.Range("F32").Value = 105
$ws.Range("F33").Value = 800
$ws.Range("F35").Value = 754
$ws.Range("F37").Value = 843
$ws.Range("F38").Value = 147556
$ws.Range("F39").Value = 836
$ws.Range("F41").Value = 1286
$ws.Range("F42").Value = 845
$ws.Range("F43").Value = 773
$ws.Range("F44").Value = 1433
$ws.Range("F45").Value = 60
$ws.Range("F46").Value = 739
$ws.Range("F47").Value = 82
$ws.Range("F48").Value = 727

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F17").Value = 135
$ws.Range("F18").Value = 1208
$ws.Range("F19").Value = 48
$ws.Range("F20").Value = 2344
$ws.Range("F21").Value = 1132
$ws.Range("F22").Value = 356
$ws.Range("F23").Value = 702
$ws.Range("F25").Value = 19
$ws.Range("F26").Value = 40
$ws.Range("F29").Value = 396
$ws.Range("F32").Value = 233
$ws.Range("F36").Value = 206
$ws.Range("F43").Value = 18
$ws.Range("F46").Value = 94

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 841
$ws.Range("F5").Value = 221
$ws.Range("F6").Value = 2600
$ws.Range("F7").Value = 4307
$ws.Range("F8").Value = 87
$ws.Range("F10").Value = 439
$ws.Range("F11").Value = 433
$ws.Range("F12").Value = 325
$ws.Range("F13").Value = 284
$ws.Range("F14").Value = 107

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1788
$ws.Range("F4").Value = 841
$ws.Range("F5").Value = 4307
$ws.Range("F6").Value = 624
$ws.Range("F7").Value = 433
$ws.Range("F8").Value = 1761
$ws.Range("F11").Value = 284
$ws.Range("F12").Value = 284
$ws.Range("F13").Value = 107
$ws.Range("F15").Value = 1042
$ws.Range("F17").Value = 1209
$ws.Range("F18").Value = 465
$ws.Range("F19").Value = 465
$ws.Range("F21").Value = 379
$ws.Range("F22").Value = 93
$ws.Range("F23").Value = 2344
$ws.Range("F24").Value = 1132
$ws.Range("F25").Value = 1160
$ws.Range("F26").Value = 1138
$ws.Range("F27").Value = 1270
$ws.Range("F29").Value = 40
$ws.Range("F30").Value = 1454
$ws.Range("F31").Value = 737
$ws.Range("F32").Value = 396
$ws.Range("F33").Value = 800
$ws.Range("F34").Value = 754
$ws.Range("F36").Value = 843
$ws.Range("F37").Value = 233
$ws.Range("F38").Value = 836
$ws.Range("F40").Value = 845
$ws.Range("F42").Value = 773
$ws.Range("F44").Value = 1433
$ws.Range("F45").Value = 60
$ws.Range("F48").Value = 739
$ws.Range("F49").Value = 727
$ws.Range("F51").Value = 94
